$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = "Datos actualizados a 26 de Marzo de 2020 a las 14:12"
$ws.Cells.Item(6, 1).Value = "Estados Unidos"
$row6 = @(68774, 563, 428, 67309, 1455, 10, 1037)
for ($i = 0; $i -lt $row6.Length; $i++) { $ws.Cells.Item(6, 2 + $i).Value = $row6[$i] }

$ws.Cells.Item(11, 1).Value = "Suiza"
$row11 = @(11316, 419, 131, 11014, 141, 18, 171)
for ($i = 0; $i -lt $row11.Length; $i++) { $ws.Cells.Item(11, 2 + $i).Value = $row11[$i] }

$ws.Cells.Item(14, 1).Value = "Paises Bajos"
$row14 = @(7431, 1019, 3, 6994, 582, 78, 434)
for ($i = 0; $i -lt $row14.Length; $i++) { $ws.Cells.Item(14, 2 + $i).Value = $row14[$i] }

$ws.Cells.Item(18, 1).Value = "Canada"
$row18 = @(3409, 0, 185, 3188, 35, 0, 36)
for ($i = 0; $i -lt $row18.Length; $i++) { $ws.Cells.Item(18, 2 + $i).Value = $row18[$i] }

$ws.Cells.Item(19, 1).Value = "Noruega"
$row19 = @(3246, 162, 6, 3226, 70, 0, 14)
for ($i = 0; $i -lt $row19.Length; $i++) { $ws.Cells.Item(19, 2 + $i).Value = $row19[$i] }

$ws.Cells.Item(20, 1).Value = "Australia"
$row20 = @(2806, 130, 170, 2623, 11, 2, 13)
for ($i = 0; $i -lt $row20.Length; $i++) { $ws.Cells.Item(20, 2 + $i).Value = $row20[$i] }

$ws.Cells.Item(34, 1).Value = "Polonia"
$row34 = @(1120, 69, 7, 1099, 3, 0, 14)
for ($i = 0; $i -lt $row34.Length; $i++) { $ws.Cells.Item(34, 2 + $i).Value = $row34[$i] }

$ws.Cells.Item(36, 1).Value = "Rumania"
$row36 = @(1029, 123, 94, 917, 29, 1, 18)
for ($i = 0; $i -lt $row36.Length; $i++) { $ws.Cells.Item(36, 2 + $i).Value = $row36[$i] }

$ws.Cells.Item(37, 1).Value = "Arabia Saudita"
$row37 = @(1012, 112, 33, 976, 6, 1, 3)
for ($i = 0; $i -lt $row37.Length; $i++) { $ws.Cells.Item(37, 2 + $i).Value = $row37[$i] }

$ws.Cells.Item(38, 1).Value = "Finlandia"
$row38 = @(928, 48, 10, 913, 22, 2, 5)
for ($i = 0; $i -lt $row38.Length; $i++) { $ws.Cells.Item(38, 2 + $i).Value = $row38[$i] }

$ws.Cells.Item(42, 1).Value = "Islandia"
$row42 = @(802, 65, 68, 732, 11, 0, 2)
for ($i = 0; $i -lt $row42.Length; $i++) { $ws.Cells.Item(42, 2 + $i).Value = $row42[$i] }

$ws.Cells.Item(43, 1).Value = "India"
$row43 = @(716, 59, 45, 657, 0, 2, 14)
for ($i = 0; $i -lt $row43.Length; $i++) { $ws.Cells.Item(43, 2 + $i).Value = $row43[$i] }

$ws.Cells.Item(44, 1).Value = "Crucero"
$row44 = @(712, 0, 597, 105, 15, 0, 10)
for ($i = 0; $i -lt $row44.Length; $i++) { $ws.Cells.Item(44, 2 + $i).Value = $row44[$i] }

$ws.Cells.Item(45, 1).Value = "Sudafrica"
$row45 = @(709, 0, 12, 697, 2, 0, 0)
for ($i = 0; $i -lt $row45.Length; $i++) { $ws.Cells.Item(45, 2 + $i).Value = $row45[$i] }

$ws.Cells.Item(46, 1).Value = "Filipinas"
$row46 = @(707, 71, 28, 634, 1, 7, 45)
for ($i = 0; $i -lt $row46.Length; $i++) { $ws.Cells.Item(46, 2 + $i).Value = $row46[$i] }

$ws.Cells.Item(75, 1).Value = "Principado de Andorra"
$row75 = @(224, 36, 1, 220, 6, 2, 3)
for ($i = 0; $i -lt $row75.Length; $i++) { $ws.Cells.Item(75, 2 + $i).Value = $row75[$i] }

$ws.Cells.Item(76, 1).Value = "Uruguay"
$row76 = @(217, 0, 0, 217, 3, 0, 0)
for ($i = 0; $i -lt $row76.Length; $i++) { $ws.Cells.Item(76, 2 + $i).Value = $row76[$i] }

$ws.Cells.Item(77, 1).Value = "San Marino"
$row77 = @(208, 0, 4, 183, 12, 0, 21)
for ($i = 0; $i -lt $row77.Length; $i++) { $ws.Cells.Item(77, 2 + $i).Value = $row77[$i] }

$ws.Cells.Item(78, 1).Value = "Kuwait"
$row78 = @(208, 13, 49, 159, 7, 0, 0)
for ($i = 0; $i -lt $row78.Length; $i++) { $ws.Cells.Item(78, 2 + $i).Value = $row78[$i] }

$ws.Cells.Item(79, 1).Value = "Costa Rica"
$row79 = @(201, 0, 2, 197, 4, 0, 2)
for ($i = 0; $i -lt $row79.Length; $i++) { $ws.Cells.Item(79, 2 + $i).Value = $row79[$i] }

$ws.Cells.Item(80, 1).Value = "Republica de Macedonia"
$row80 = @(201, 24, 3, 195, 1, 0, 3)
for ($i = 0; $i -lt $row80.Length; $i++) { $ws.Cells.Item(80, 2 + $i).Value = $row80[$i] }

$ws.Cells.Item(128, 1).Value = "Kenia"
$row128 = @(31, 3, 1, 30, 0, 0, 0)
for ($i = 0; $i -lt $row128.Length; $i++) { $ws.Cells.Item(128, 2 + $i).Value = $row128[$i] }

$ws.Cells.Item(129, 1).Value = "Macao"
$row129 = @(31, 0, 10, 21, 0, 0, 0)
for ($i = 0; $i -lt $row129.Length; $i++) { $ws.Cells.Item(129, 2 + $i).Value = $row129[$i] }

$ws.Cells.Item(133, 1).Value = "Isla de Man"
$row133 = @(25, 2, 0, 25, 0, 0, 0)
for ($i = 0; $i -lt $row133.Length; $i++) { $ws.Cells.Item(133, 2 + $i).Value = $row133[$i] }

$ws.Cells.Item(134, 1).Value = "Polinesia Francesa"
$row134 = @(25, 0, 0, 25, 0, 0, 0)
for ($i = 0; $i -lt $row134.Length; $i++) { $ws.Cells.Item(134, 2 + $i).Value = $row134[$i] }

$ws.Cells.Item(136, 1).Value = "Madagascar"
$row136 = @(23, 4, 0, 23, 0, 0, 0)
for ($i = 0; $i -lt $row136.Length; $i++) { $ws.Cells.Item(136, 2 + $i).Value = $row136[$i] }

$ws.Cells.Item(137, 1).Value = "Togo"
$row137 = @(23, 0, 1, 22, 0, 0, 0)
for ($i = 0; $i -lt $row137.Length; $i++) { $ws.Cells.Item(137, 2 + $i).Value = $row137[$i] }
